# Insert a new data row at row 79 (shifts existing rows 79..202 down to 80..203)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("79:79").Insert()

$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44540
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112037
$ws.Range("G79").Value = "Cebollín"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 180
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 6000
$ws.Range("M79").Value = 6000
$ws.Range("N79").Value = "$/paquete 36 unidades"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 167
$ws.Range("Q79").Value = 36
$ws.Range("R79").Value = "Hortaliza"
